$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix Jhonny Bond's phone number
$ws.Range("D4").Value = "+336748195"

# Update the header "Name" -> "Nom"
$ws.Range("B1").Value = "Nom"

# Reselect cell C1, matching the saved view state in the workbook
$ws.Range("C1").Select()
